$d = $word.ActiveDocument
$apos = [char]0x2019

# -----------------------------------------------------------------------
# Hunk 1: "...doesn't improve delays but rather switches..." gains a
# parenthetical clarification in the middle, splitting the original run
# into three runs with identical formatting.
# -----------------------------------------------------------------------
$target1 = "decreasing discount factor or increasing reward doesn${apos}t improve delays but rather switches between delay or no delay at all."
$range1 = $d.Content
$range1.Find.Execute($target1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base1 = $range1.Start

$h1a = "decreasing discount factor or increasing reward doesn${apos}t improve delays"
$h1b = " (much, does happen a little bit in some param settings)"

# Pin the boundary right after the sentence so the untouched, identically
# formatted run that immediately follows it (a lone space) is not swept
# up into the freshly-edited run when the insertion below triggers a
# re-merge of adjacent same-format runs.
$endOfTarget1 = $base1 + $target1.Length
$protectEnd1 = $d.Range($endOfTarget1, $endOfTarget1)
$d.Bookmarks.Add("ZZ_PROTECT_1_END", $protectEnd1) | Out-Null

$ins1 = $d.Range($base1 + $h1a.Length, $base1 + $h1a.Length)
$ins1.InsertAfter($h1b)
# Toggling a character property on the just-inserted text (and back to its
# original value) forces it to stay a distinct run instead of re-merging
# with its identically-formatted neighbours.
$ins1.Bold = $true
$ins1.Bold = $false

$d.Bookmarks("ZZ_PROTECT_1_END").Delete()

# -----------------------------------------------------------------------
# Hunk 2: "...then? What about hyperbolic discounting?" gains a
# parenthetical "(atleast for all patterns)" in the middle, splitting the
# original run into five runs, and a _GoBack bookmark is left right after
# the resulting text.
# -----------------------------------------------------------------------
$target2 = ". This is probably not the best explanation then? What about hyperbolic discounting?"
$range2 = $d.Content
$range2.Find.Execute($target2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base2 = $range2.Start

$h2a = ". This is probably not the best explanation then"
$h2b = " ("
$h2c = "atleast"
$h2d = " for all patterns)"
$h2e = "? What about hyperbolic discounting?"

# Pin the boundary right before this sentence too, so the untouched,
# identically formatted run immediately preceding it ("... only small
# gaps") does not get absorbed by the inserts below.
$protectStart2 = $d.Range($base2, $base2)
$d.Bookmarks.Add("ZZ_PROTECT_2_START", $protectStart2) | Out-Null

# Insert each new fragment in turn. Because later inserts can re-coalesce
# earlier splits that share identical formatting, each split point is
# pinned in place with a temporary bookmark until every fragment has been
# inserted; the temporary bookmarks are then removed as the very last
# step so the run boundaries they protected are left standing.
$off2 = $base2 + $h2a.Length
$insB = $d.Range($off2, $off2)
$insB.InsertAfter($h2b)
$d.Bookmarks.Add("ZZ_SPLIT_B", $insB) | Out-Null
$off2 = $off2 + $h2b.Length

$insC = $d.Range($off2, $off2)
$insC.InsertAfter($h2c)
$d.Bookmarks.Add("ZZ_SPLIT_C", $insC) | Out-Null
$off2 = $off2 + $h2c.Length

$insD = $d.Range($off2, $off2)
$insD.InsertAfter($h2d)
$d.Bookmarks.Add("ZZ_SPLIT_D", $insD) | Out-Null
$off2 = $off2 + $h2d.Length

$d.Bookmarks("ZZ_PROTECT_2_START").Delete()
$d.Bookmarks("ZZ_SPLIT_B").Delete()
$d.Bookmarks("ZZ_SPLIT_C").Delete()
$d.Bookmarks("ZZ_SPLIT_D").Delete()

$endOfHunk2 = $base2 + $h2a.Length + $h2b.Length + $h2c.Length + $h2d.Length + $h2e.Length

# -----------------------------------------------------------------------
# Move the _GoBack bookmark from its old spot (after "commitments" near
# the end of the document) to the collapsed position right after the
# text that was just edited above.
# -----------------------------------------------------------------------
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

$newGoBackRange = $d.Range($endOfHunk2, $endOfHunk2)
$d.Bookmarks.Add("_GoBack", $newGoBackRange) | Out-Null
